# Weekly crime data refresh: shift the reporting week forward by one week
# (Volume/Number bump, date range bump) and update the underlying crime
# statistics table (rows 16, 17, 19, 21, 24, 25) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings) - Volume/Number and date range
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/20/2023  Through  11/26/2023"

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 20
$ws.Range("K16").Value = -28.571428571428
$ws.Range("L16").Value = 42.857142857142
$ws.Range("M16").Value = -13.043478260869
$ws.Range("N16").Value = -89.637305699481

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("L17").Value = -20

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny (C19 switches from numeric 1 to text "0", and the
# old numeric value 2 moves from F19 down to 1)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = "'0"
$ws.Range("D20").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("F19").Value = 1
$ws.Range("H19").Value = -50
$ws.Range("L19").Value = 84
$ws.Range("N19").Value = -73.255813953488

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 1
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = -40
$ws.Range("I21").Value = 81
$ws.Range("K21").Value = 12.5
$ws.Range("L21").Value = 44.642857142857
$ws.Range("M21").Value = -17.346938775510
$ws.Range("N21").Value = -82.275711159737

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny (C24 switches from text "0" to numeric 1)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 1
$ws.Range("F24").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 400
$ws.Range("I24").Value = 39
$ws.Range("K24").Value = 34.482758620689
$ws.Range("L24").Value = 21.875
$ws.Range("M24").Value = -59.375

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault (C25/D25 switch from numeric to text "0", E25
# switches from numeric 100 to text "***.*"; old F25 value 3 -> 2)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = "'0"
$ws.Range("D20").Copy()
$ws.Range("C25").PasteSpecial(-4122)

$ws.Range("D25").Value = "'0"
$ws.Range("D20").Copy()
$ws.Range("D25").PasteSpecial(-4122)

$ws.Range("E25").Value = "***.*"
$ws.Range("E20").Copy()
$ws.Range("E25").PasteSpecial(-4122)

$ws.Range("F25").Value = 2
$ws.Range("H25").Value = 0
$ws.Range("L25").Value = 95.833333333333
